$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").Font.Name = "Consolas"
$ws.Range("D2").Font.Name = "Courier New"
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D4").Font.Name = "Times New Roman"
$ws.Range("D5").Font.Name = "Comic Sans MS"
